# Word COM-interop edit script
# Target paragraph (before):
#   "The security model implemented over the directed graph is Biba-Bell
#    Lapadulla [6]. Biba is a model focused on integrity and Bell Lapudalla
#    is focused on confidentiality. These models describe information flow
#    between different subjects."
# Target paragraph (after):
#   "The security models implemented over the directed graph are Biba and
#    Bell-LaPadula [6]. Biba is a model focused on integrity and
#    Bell-LaPadula is focused on confidentiality. These models describe
#    information flow between different subjects."

$d = $word.ActiveDocument

# 1) "model" -> "models" (pluralize), keep surrounding text intact.
$d.Content.Find.Execute("The security model implemented", $true, $false, $false, $false, $false, $true, 1, $false, "The security models implemented", 2)

# 2) "is" -> "are" right before the model names.
$d.Content.Find.Execute("directed graph is Biba", $true, $false, $false, $false, $false, $true, 1, $false, "directed graph are Biba", 2)

# 3) "Biba-Bell Lapadulla" -> "Biba and Bell-LaPadula" (fix the misspelled
#    "Lapadulla" to "LaPadula" and change the model-name conjunction).
$d.Content.Find.Execute("Biba-Bell Lapadulla [6].", $true, $false, $false, $false, $false, $true, 1, $false, "Biba and Bell-LaPadula [6].", 2)

# 4) "Bell Lapudalla" -> "Bell-LaPadula" (second mention, also misspelled).
$d.Content.Find.Execute("Bell Lapudalla is focused", $true, $false, $false, $false, $false, $true, 1, $false, "Bell-LaPadula is focused", 2)
